# edit.ps1
# Applies the "cryptos list" refresh described by the commit:
#   "Updated cryptos list on Fri Jun 30 11:28:09 UTC 2023 with GitHub Actions"
#
# The sheet (Sheet1) holds a table with columns:
#   A = rank index (untouched), B = Coin name, C = Link, D = Price, E = Volume(1h)
#
# Many values in column D look like plain numbers (e.g. "1.0000", "0.000007621")
# but must stay as literal TEXT (they are thousand/European-style formatted price
# strings, not numeric cells). Setting NumberFormat="@" (Text) right before writing
# the value prevents Excel from re-parsing/normalizing the string into a number,
# and resetting .Style back to "Normal" afterwards keeps the cell formatting exactly
# as it was (no explicit style index left behind).
function Set-TextCell {
    param($ws, $ref, $val)
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
Set-TextCell $ws 'D2' '30.803.82'
Set-TextCell $ws 'E2' '  +0.61%  '

# Row 3
Set-TextCell $ws 'D3' '1.890.15'
Set-TextCell $ws 'E3' '  +1.16%  '

# Row 4
Set-TextCell $ws 'D4' '1.0000'
Set-TextCell $ws 'E4' '  -0.01%  '

# Row 5
Set-TextCell $ws 'D5' '240.00'
Set-TextCell $ws 'E5' '  +2.03%  '

# Row 6
Set-TextCell $ws 'E6' '  +0.03%  '

# Row 7
Set-TextCell $ws 'D7' '0.4803'
Set-TextCell $ws 'E7' '  +1.99%  '

# Row 8
Set-TextCell $ws 'D8' '0.2958'
Set-TextCell $ws 'E8' '  +6.58%  '

# Row 9
Set-TextCell $ws 'D9' '0.06624'
Set-TextCell $ws 'E9' '  +3.80%  '

# Row 10
Set-TextCell $ws 'D10' '18.75'
Set-TextCell $ws 'E10' '  +3.13%  '

# Row 11
Set-TextCell $ws 'D11' '101.33'
Set-TextCell $ws 'E11' '  +18.71%  '

# Row 12
Set-TextCell $ws 'D12' '1.885.42'
Set-TextCell $ws 'E12' '  +0.95%  '

# Row 13
Set-TextCell $ws 'D13' '0.07641'
Set-TextCell $ws 'E13' '  +2.37%  '

# Row 14
Set-TextCell $ws 'D14' '5.144'
Set-TextCell $ws 'E14' '  +3.25%  '

# Row 15
Set-TextCell $ws 'D15' '0.6585'
Set-TextCell $ws 'E15' '  +3.47%  '

# Row 16
Set-TextCell $ws 'D16' '307.68'
Set-TextCell $ws 'E16' '  +27.30%  '

# Row 17
Set-TextCell $ws 'D17' '30.779.62'
Set-TextCell $ws 'E17' '  +0.65%  '

# Row 18
Set-TextCell $ws 'E18' '  +2.89%  '

# Row 19/20: ShibaInu and Dai swap positions (full row content exchange)
# Row 19
Set-TextCell $ws 'B19' 'Dai'
Set-TextCell $ws 'C19' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextCell $ws 'D19' '1.000'
Set-TextCell $ws 'E19' '  +0.09%  '

# Row 20
Set-TextCell $ws 'B20' 'ShibaInu'
Set-TextCell $ws 'C20' 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextCell $ws 'D20' '0.000007621'
Set-TextCell $ws 'E20' '  +3.18%  '

# Row 21
Set-TextCell $ws 'D21' '2.137.35'
Set-TextCell $ws 'E21' '  +1.89%  '

# Row 22
Set-TextCell $ws 'D22' '1.000'
Set-TextCell $ws 'E22' '  -0.07%  '

# Row 23
Set-TextCell $ws 'E23' '  +3.28%  '

# Row 24
Set-TextCell $ws 'D24' '6.188'
Set-TextCell $ws 'E24' '  +2.76%  '

# Row 25
Set-TextCell $ws 'E25' '  -0.81%  '

# Row 26
Set-TextCell $ws 'D26' '167.93'
Set-TextCell $ws 'E26' '  +1.10%  '

# Row 27
Set-TextCell $ws 'D27' '20.57'
Set-TextCell $ws 'E27' '  +12.70%  '

# Row 28
Set-TextCell $ws 'D28' '1.961'
Set-TextCell $ws 'E28' '  +3.59%  '

# Row 29
Set-TextCell $ws 'D29' '0.1129'
Set-TextCell $ws 'E29' '  +9.69%  '

# Row 30
Set-TextCell $ws 'D30' '1.346'
Set-TextCell $ws 'E30' '  -2.56%  '

# Row 31
Set-TextCell $ws 'D31' '4.193'
Set-TextCell $ws 'E31' '  +2.13%  '

# Row 32
Set-TextCell $ws 'D32' '4.012'
Set-TextCell $ws 'E32' '  +3.53%  '

# Row 33
Set-TextCell $ws 'D33' '0.05105'
Set-TextCell $ws 'E33' '  +3.33%  '

# Row 34
Set-TextCell $ws 'D34' '0.7435'
Set-TextCell $ws 'E34' '  +4.86%  '

# Row 35
Set-TextCell $ws 'D35' '1.166'
Set-TextCell $ws 'E35' '  +0.84%  '

# Row 36
Set-TextCell $ws 'D36' '2.717'
Set-TextCell $ws 'E36' '  +0.53%  '

# Row 37
Set-TextCell $ws 'D37' '0.01980'
Set-TextCell $ws 'E37' '  +4.08%  '

# Row 38
Set-TextCell $ws 'D38' '2.708'
Set-TextCell $ws 'E38' '  +0.49%  '

# Row 39
Set-TextCell $ws 'D39' '2.060'
Set-TextCell $ws 'E39' '  +3.17%  '

# Row 40/41: TrustWalletToken and Quant swap positions (full row content exchange)
# Row 40
Set-TextCell $ws 'B40' 'Quant'
Set-TextCell $ws 'C40' 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextCell $ws 'D40' '109.56'
Set-TextCell $ws 'E40' '  +3.40%  '

# Row 41
Set-TextCell $ws 'B41' 'TrustWalletToken'
Set-TextCell $ws 'C41' 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextCell $ws 'D41' '0.8762'
Set-TextCell $ws 'E41' '  -0.53%  '

# Row 42
Set-TextCell $ws 'E42' '  +0.03%  '

# Row 43
Set-TextCell $ws 'D43' '0.4217'
Set-TextCell $ws 'E43' '  +2.45%  '

# Row 44
Set-TextCell $ws 'E44' '  +2.07%  '

# Row 45
Set-TextCell $ws 'D45' '67.76'
Set-TextCell $ws 'E45' '  +8.74%  '

# Row 46
Set-TextCell $ws 'D46' '7.386'
Set-TextCell $ws 'E46' '  -1.15%  '

# Row 47
Set-TextCell $ws 'D47' '9.148'
Set-TextCell $ws 'E47' '  +5.56%  '

# Row 48
Set-TextCell $ws 'E48' '  +0.39%  '

# Row 49
Set-TextCell $ws 'D49' '34.99'
Set-TextCell $ws 'E49' '  +3.85%  '

# Row 50
Set-TextCell $ws 'D50' '0.05654'
Set-TextCell $ws 'E50' '  +1.30%  '

# Row 51
Set-TextCell $ws 'D51' '1.403'
Set-TextCell $ws 'E51' '  +1.00%  '

Write-Output "Updated cryptos list: applied $([int]101) cell changes across 50 data rows"
